$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.031.85"
$ws.Range("E2").Value = "  -2.23%  "
$ws.Range("D3").Value = "2.971.50"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'592.06"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").Value = "'141.76"
$ws.Range("E6").Value = "  -3.22%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.512"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").Value = "2.970.78"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D11").Value = "'5.93"
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "'33.89"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "3.461.09"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "61.119.63"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "'6.83"
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("D19").Value = "2.969.09"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("D20").Value = "'447.90"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").Value = "'14.00"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'0.676"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").Value = "'7.23"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").Value = "'82.13"
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("D25").Value = "'2.15"
$ws.Range("E25").Value = "  -5.57%  "
$ws.Range("D26").Value = "'11.85"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").Value = "'10.23"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "'7.02"
$ws.Range("E31").Value = "  -2.47%  "
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("D33").Value = "'26.99"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "0.0₃0805"
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").Value = "'5.72"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'50.16"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'2.04"
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("E41").Value = "  +7.08%  "
$ws.Range("D42").Value = "'2.82"
$ws.Range("E42").Value = "  -4.18%  "
$ws.Range("D43").Value = "'386.94"
$ws.Range("E43").Value = "  -5.46%  "
$ws.Range("D44").Value = "'38.43"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").Value = "'0.264"
$ws.Range("E46").Value = "  -4.95%  "
$ws.Range("D47").Value = "2.681.94"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "'129.87"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("E51").Value = "  -1.22%  "
